# "rozmery po uprave velikosti boxu" - add rows for the resized box
# (device / orientation / resolution / box size / letter size) and
# switch the sheet to portrait A4 for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Row 11 references the new shared string "Nexus 6P" - write it first so
# it lands before "DP" (row 10) in the shared-string table, matching the
# order the strings were first used in the original edit.
$ws.Range("A11").Value = "Nexus 6P"

# Row 10: DP - box only, no device resolution
$ws.Range("A10").Value = "DP"
$ws.Range("E10").Value = 120
$ws.Range("F10").Value = 180
$ws.Range("J10").Formula = "=E10/F10"

# Row 11: Nexus 6P device resolution
$ws.Range("B11").Value = "landscape"
$ws.Range("C11").Value = 2560
$ws.Range("D11").Value = 1440

# Row 12: LeliMath box, landscape
$ws.Range("A12").Value = "LeliMath"
$ws.Range("B12").Value = "landscape"
$ws.Range("C12").Value = 2390
$ws.Range("D12").Value = 1190
$ws.Range("E12").Value = 420
$ws.Range("F12").Value = 630
$ws.Range("G12").Value = 29
$ws.Range("H12").Value = 37
$ws.Range("J12").Formula = "=E12/F12"

# Row 13: Nexus 5x device resolution
$ws.Range("A13").Value = "Nexus 5x"
$ws.Range("C13").Value = 1920
$ws.Range("D13").Value = 1080

# Row 14: LeliMath box, portrait
$ws.Range("A14").Value = "LeliMath"
$ws.Range("B14").Value = "portrait"
$ws.Range("C14").Value = 1080
$ws.Range("D14").Value = 1590
$ws.Range("E14").Value = 310
$ws.Range("F14").Value = 470
$ws.Range("G14").Value = 24
$ws.Range("H14").Value = 30
$ws.Range("J14").Formula = "=E14/F14"

# Match the author's final cursor position
$ws.Range("J15").Select()

# Page setup: A4, portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
